# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates per the commit diff, sheet by sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3034635.5
$ws.Range("J17").Value = 3337990.2
$ws.Range("L17").Value = 10013970.6
$ws.Range("N17").Value = -10014306.6
$ws.Range("H64").Value = 5300
$ws.Range("J64").Value = 5300
$ws.Range("L64").Value = 5300
$ws.Range("N64").Value = -5796
$ws.Range("H67").Value = 5300
$ws.Range("J67").Value = 5300
$ws.Range("L67").Value = 5300
$ws.Range("N67").Value = -7016
$ws.Range("H131").Value = 1236.875
$ws.Range("I131").Value = 1236.875
$ws.Range("K131").Value = 3710.625
$ws.Range("M131").Value = 1329.375
$ws.Range("H137").Value = 41424.28
$ws.Range("I137").Value = 1120.4
$ws.Range("K137").Value = 3361.2
$ws.Range("M137").Value = -811.2000000000003
$ws.Range("H138").Value = 13160509
$ws.Range("I138").Value = 37038280
$ws.Range("J138").Value = 3369.0815
$ws.Range("K138").Value = 111114840
$ws.Range("L138").Value = 10107.2445
$ws.Range("M138").Value = -111109700
$ws.Range("N138").Value = -20387.2445
$ws.Range("H141").Value = 1105.0638
$ws.Range("I141").Value = 805.5714
$ws.Range("J141").Value = 3620.8
$ws.Range("K141").Value = 2416.7142
$ws.Range("L141").Value = 10862.4
$ws.Range("M141").Value = 2763.2858
$ws.Range("N141").Value = -21222.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1577.2142
$ws.Range("I2").Value = 1398.1
$ws.Range("J2").Value = 2025
$ws.Range("K2").Value = 1398.1
$ws.Range("L2").Value = 2025
$ws.Range("M2").Value = -1285.1
$ws.Range("N2").Value = -2251
$ws.Range("H32").Value = 3539.43
$ws.Range("I32").Value = 3389.1135
$ws.Range("K32").Value = 3389.1135
$ws.Range("M32").Value = -3102.1135
$ws.Range("H61").Value = 784046.9399999999
$ws.Range("I61").Value = 1001572.56
$ws.Range("J61").Value = 954.8
$ws.Range("K61").Value = 1001572.56
$ws.Range("L61").Value = 954.8
$ws.Range("M61").Value = -1001360.56
$ws.Range("N61").Value = -1378.8
$ws.Range("H116").Value = 1577.2142
$ws.Range("I116").Value = 1398.1
$ws.Range("J116").Value = 2025
$ws.Range("K116").Value = 1398.1
$ws.Range("L116").Value = 2025
$ws.Range("M116").Value = 895.9000000000001
$ws.Range("N116").Value = -6613
$ws.Range("H122").Value = 1788.1471
$ws.Range("I122").Value = 1371.8966
$ws.Range("J122").Value = 4202.4
$ws.Range("K122").Value = 4115.6898
$ws.Range("L122").Value = 12607.2
$ws.Range("M122").Value = -1665.6898
$ws.Range("N122").Value = -17507.2
$ws.Range("H132").Value = 10102.946
$ws.Range("I132").Value = 1085.5
$ws.Range("J132").Value = 85248.336
$ws.Range("K132").Value = 3256.5
$ws.Range("L132").Value = 255745.008
$ws.Range("M132").Value = -726.5
$ws.Range("N132").Value = -260805.008
$ws.Range("H136").Value = 784046.9399999999
$ws.Range("I136").Value = 1001572.56
$ws.Range("J136").Value = 954.8
$ws.Range("K136").Value = 3004717.68
$ws.Range("L136").Value = 2864.4
$ws.Range("M136").Value = -3002167.68
$ws.Range("N136").Value = -7964.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1577.2142
$ws.Range("I3").Value = 1398.1
$ws.Range("J3").Value = 2025
$ws.Range("K3").Value = 1398.1
$ws.Range("L3").Value = 2025
$ws.Range("M3").Value = -1284.1
$ws.Range("N3").Value = -2253
$ws.Range("H57").Value = 34775
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 34775
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 34775
$ws.Range("N57").Value = -36215
$ws.Range("M57").ClearContents()
$ws.Range("H105").Value = 2908.0625
$ws.Range("I105").Value = 2911.7273
$ws.Range("J105").Value = 2900
$ws.Range("K105").Value = 2911.7273
$ws.Range("L105").Value = 2900
$ws.Range("M105").Value = -1164.7273
$ws.Range("N105").Value = -6394
$ws.Range("H134").Value = 23736.166
$ws.Range("I134").Value = 27356.635
$ws.Range("K134").Value = 82069.905
$ws.Range("M134").Value = -79534.905
$ws.Range("H136").Value = 34775
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 34775
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 34775
$ws.Range("N136").Value = -44975
$ws.Range("M136").ClearContents()
$ws.Range("H140").Value = 46280
$ws.Range("J140").Value = 46280
$ws.Range("L140").Value = 46280
$ws.Range("N140").Value = -56640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2991.535
$ws.Range("I31").Value = 1772.25
$ws.Range("J31").Value = 5267.533
$ws.Range("K31").Value = 1772.25
$ws.Range("L31").Value = 5267.533
$ws.Range("M31").Value = -1477.25
$ws.Range("N31").Value = -5857.533
$ws.Range("H34").Value = 2991.535
$ws.Range("I34").Value = 1772.25
$ws.Range("J34").Value = 5267.533
$ws.Range("K34").Value = 1772.25
$ws.Range("L34").Value = 5267.533
$ws.Range("M34").Value = -1570.25
$ws.Range("N34").Value = -5671.533
$ws.Range("H62").Value = 76926696
$ws.Range("I62").Value = 90911910
$ws.Range("J62").Value = 8003
$ws.Range("K62").Value = 90911910
$ws.Range("L62").Value = 8003
$ws.Range("M62").Value = -90911286
$ws.Range("N62").Value = -9251
$ws.Range("H65").Value = 76926696
$ws.Range("I65").Value = 90911910
$ws.Range("J65").Value = 8003
$ws.Range("K65").Value = 454559550
$ws.Range("L65").Value = 40015
$ws.Range("M65").Value = -454556430
$ws.Range("N65").Value = -46255
$ws.Range("H122").Value = 4000.2
$ws.Range("I122").Value = 5333.6665
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 16000.9995
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -13550.9995
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 1705.6041
$ws.Range("I132").Value = 1336.4667
$ws.Range("J132").Value = 7242.6665
$ws.Range("K132").Value = 4009.4001
$ws.Range("L132").Value = 21727.9995
$ws.Range("M132").Value = -1479.4001
$ws.Range("N132").Value = -26787.9995
$ws.Range("H134").Value = 741.62
$ws.Range("I134").Value = 674.1778
$ws.Range("J134").Value = 1348.6
$ws.Range("K134").Value = 2022.5334
$ws.Range("L134").Value = 4045.8
$ws.Range("M134").Value = 512.4665999999997
$ws.Range("N134").Value = -9115.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 297.6
$ws.Range("I8").Value = 297.6
$ws.Range("K8").Value = 892.8000000000001
$ws.Range("M8").Value = -753.8000000000001
$ws.Range("H102").Value = 7000
$ws.Range("J102").Value = 7000
$ws.Range("L102").Value = 21000
$ws.Range("N102").Value = -25868
$ws.Range("H131").Value = 696.67
$ws.Range("J131").Value = 713.7717
$ws.Range("L131").Value = 2141.3151
$ws.Range("N131").Value = -12221.3151
$ws.Range("H140").Value = 2916.25
$ws.Range("I140").Value = 2599.8
$ws.Range("J140").Value = 3142.2856
$ws.Range("K140").Value = 7799.400000000001
$ws.Range("L140").Value = 9426.856800000001
$ws.Range("M140").Value = -2619.400000000001
$ws.Range("N140").Value = -19786.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3866.4443
$ws.Range("I113").Value = 3100
$ws.Range("J113").Value = 4479.6
$ws.Range("K113").Value = 3100
$ws.Range("L113").Value = 4479.6
$ws.Range("M113").Value = -930
$ws.Range("N113").Value = -8819.6
$ws.Range("H122").Value = 57971796
$ws.Range("I122").Value = 16667381
$ws.Range("J122").Value = 333334560
$ws.Range("K122").Value = 50002143
$ws.Range("L122").Value = 1000003680
$ws.Range("M122").Value = -49999693
$ws.Range("N122").Value = -1000008580
$ws.Range("H132").Value = 87741.664
$ws.Range("I132").Value = 63491.176
$ws.Range("K132").Value = 190473.528
$ws.Range("M132").Value = -187943.528

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 224.76923
$ws.Range("I16").Value = 224.76923
$ws.Range("K16").Value = 224.76923
$ws.Range("M16").Value = -54.76922999999999
$ws.Range("H22").Value = 1753.9231
$ws.Range("I22").Value = 1242.8
$ws.Range("J22").Value = 3457.6667
$ws.Range("K22").Value = 1242.8
$ws.Range("L22").Value = 3457.6667
$ws.Range("M22").Value = -947.8
$ws.Range("N22").Value = -4047.6667
$ws.Range("H27").Value = 1753.9231
$ws.Range("I27").Value = 1242.8
$ws.Range("J27").Value = 3457.6667
$ws.Range("K27").Value = 1242.8
$ws.Range("L27").Value = 3457.6667
$ws.Range("M27").Value = -1135.8
$ws.Range("N27").Value = -3671.6667
$ws.Range("H55").Value = 197.89473
$ws.Range("I55").Value = 190.66667
$ws.Range("J55").Value = 225
$ws.Range("K55").Value = 190.66667
$ws.Range("L55").Value = 225
$ws.Range("M55").Value = -17.66667000000001
$ws.Range("N55").Value = -571
$ws.Range("H61").Value = 6399.909
$ws.Range("J61").Value = 7537.375
$ws.Range("L61").Value = 7537.375
$ws.Range("N61").Value = -7941.375
$ws.Range("H68").Value = 2356.4614
$ws.Range("I68").Value = 2427.7778
$ws.Range("J68").Value = 2196
$ws.Range("K68").Value = 2427.7778
$ws.Range("L68").Value = 2196
$ws.Range("M68").Value = -1678.7778
$ws.Range("N68").Value = -3694
$ws.Range("H71").Value = 2356.4614
$ws.Range("I71").Value = 2427.7778
$ws.Range("J71").Value = 2196
$ws.Range("K71").Value = 12138.889
$ws.Range("L71").Value = 10980
$ws.Range("M71").Value = -8394.888999999999
$ws.Range("N71").Value = -18468
$ws.Range("H93").Value = 1520.9
$ws.Range("I93").Value = 1412.1111
$ws.Range("K93").Value = 1412.1111
$ws.Range("M93").Value = -164.1111000000001
$ws.Range("H113").Value = 6399.909
$ws.Range("J113").Value = 7537.375
$ws.Range("L113").Value = 7537.375
$ws.Range("N113").Value = -11877.375
$ws.Range("H122").Value = 855277.5
$ws.Range("I122").Value = 1784916.8
$ws.Range("K122").Value = 5354750.4
$ws.Range("M122").Value = -5352300.4
$ws.Range("H132").Value = 603884.25
$ws.Range("I132").Value = 635593.9399999999
$ws.Range("K132").Value = 1906781.82
$ws.Range("M132").Value = -1904251.82
$ws.Range("H136").Value = 1057.4546
$ws.Range("I136").Value = 945.75
$ws.Range("K136").Value = 2837.25
$ws.Range("M136").Value = -287.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1931280
$ws.Range("J113").Value = 6756859
$ws.Range("L113").Value = 20270577
$ws.Range("N113").Value = -20274917
$ws.Range("H122").Value = 1162.1613
$ws.Range("I122").Value = 1170.9
$ws.Range("K122").Value = 3512.7
$ws.Range("M122").Value = -1062.7
$ws.Range("H136").Value = 19853360
